# Update "想去人数" (F column) figures across sheets to reflect the
# latest generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1026
$ws1.Range("F4").Value  = 13490
$ws1.Range("F6").Value  = 1022
$ws1.Range("F7").Value  = 16
$ws1.Range("F9").Value  = 133
$ws1.Range("F11").Value = 78
$ws1.Range("F14").Value = 13488
$ws1.Range("F16").Value = 599
$ws1.Range("F17").Value = 8948
$ws1.Range("F18").Value = 8
$ws1.Range("F19").Value = 8022
$ws1.Range("F20").Value = 250
$ws1.Range("F22").Value = 147
$ws1.Range("F26").Value = 21
$ws1.Range("F27").Value = 1020
$ws1.Range("F30").Value = 394
$ws1.Range("F32").Value = 174
$ws1.Range("F33").Value = 378
$ws1.Range("F34").Value = 94

# --- Sheet "演出" (shows) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 38

# --- Sheet "全部类型" (all types, aggregated) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1026
$ws4.Range("F4").Value  = 13490
$ws4.Range("F6").Value  = 1022
$ws4.Range("F7").Value  = 16
$ws4.Range("F9").Value  = 133
$ws4.Range("F11").Value = 78
$ws4.Range("F14").Value = 13488
$ws4.Range("F16").Value = 599
$ws4.Range("F17").Value = 8948
$ws4.Range("F18").Value = 8
$ws4.Range("F19").Value = 8022
$ws4.Range("F20").Value = 251
$ws4.Range("F22").Value = 147
$ws4.Range("F26").Value = 21
$ws4.Range("F27").Value = 1020
$ws4.Range("F30").Value = 38
$ws4.Range("F32").Value = 394
$ws4.Range("F34").Value = 174
$ws4.Range("F35").Value = 378
$ws4.Range("F36").Value = 94
